$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated win/transition probabilities for the team-specific matrix
# (changes to team matrices from games pulled march 7)
    $ws.Range("B2").Value = 0.1987767584097859
    $ws.Range("C2").Value = 0.5596330275229358
    $ws.Range("J2").Value = 0.01223241590214067
    $ws.Range("O2").Value = 0.003058103975535168
    $ws.Range("P2").Value = 0.1467889908256881
    $ws.Range("S2").Value = 0.07951070336391437
    $ws.Range("B3").Value = 0.005181347150259068
    $ws.Range("C3").Value = 0.05181347150259067
    $ws.Range("J3").Value = 0.03626943005181347
    $ws.Range("P3").Value = 0.7150259067357513
    $ws.Range("S3").Value = 0.1917098445595855
    $ws.Range("J4").Value = 0.02083333333333333
    $ws.Range("O4").Value = 0.02083333333333333
    $ws.Range("P4").Value = 0.6666666666666666
    $ws.Range("S4").Value = 0.2916666666666667
    $ws.Range("B6").Value = 0.06837606837606838
    $ws.Range("D6").Value = 0.008547008547008548
    $ws.Range("F6").Value = 0.05982905982905983
    $ws.Range("J6").Value = 0.1965811965811966
    $ws.Range("O6").Value = 0.03846153846153846
    $ws.Range("Q6").Value = 0.1581196581196581
    $ws.Range("R6").Value = 0.1025641025641026
    $ws.Range("S6").Value = 0.3675213675213675
    $ws.Range("B7").Value = 0.0962962962962963
    $ws.Range("D7").Value = 0.01851851851851852
    $ws.Range("F7").Value = 0.03333333333333333
    $ws.Range("J7").Value = 0.1740740740740741
    $ws.Range("O7").Value = 0.01111111111111111
    $ws.Range("Q7").Value = 0.1666666666666667
    $ws.Range("R7").Value = 0.08148148148148149
    $ws.Range("S7").Value = 0.4185185185185185
    $ws.Range("B8").Value = 0.07786259541984733
    $ws.Range("D8").Value = 0.01679389312977099
    $ws.Range("F8").Value = 0.07022900763358779
    $ws.Range("J8").Value = 0.1083969465648855
    $ws.Range("O8").Value = 0.01679389312977099
    $ws.Range("Q8").Value = 0.1862595419847328
    $ws.Range("R8").Value = 0.09465648854961832
    $ws.Range("S8").Value = 0.4290076335877863
    $ws.Range("B9").Value = 0.1139240506329114
    $ws.Range("D9").Value = 0.0379746835443038
    $ws.Range("F9").Value = 0.05063291139240506
    $ws.Range("J9").Value = 0.1455696202531646
    $ws.Range("O9").Value = 0.0189873417721519
    $ws.Range("Q9").Value = 0.2025316455696203
    $ws.Range("R9").Value = 0.08227848101265822
    $ws.Range("S9").Value = 0.3481012658227848
    $ws.Range("B10").Value = 0.1008064516129032
    $ws.Range("D10").Value = 0.01948924731182796
    $ws.Range("F10").Value = 0.05174731182795699
    $ws.Range("J10").Value = 0.1283602150537634
    $ws.Range("O10").Value = 0.01881720430107527
    $ws.Range("Q10").Value = 0.2345430107526882
    $ws.Range("R10").Value = 0.0793010752688172
    $ws.Range("S10").Value = 0.3669354838709677
    $ws.Range("F11").Value = 0.002398081534772182
    $ws.Range("G11").Value = 0.1702637889688249
    $ws.Range("J11").Value = 0.0815347721822542
    $ws.Range("K11").Value = 0.2014388489208633
    $ws.Range("L11").Value = 0.5275779376498801
    $ws.Range("S11").Value = 0.01678657074340528
    $ws.Range("G12").Value = 0.7685589519650655
    $ws.Range("J12").Value = 0.1528384279475982
    $ws.Range("K12").Value = 0.01746724890829694
    $ws.Range("L12").Value = 0.03056768558951965
    $ws.Range("S12").Value = 0.03056768558951965
    $ws.Range("G13").Value = 0.5882352941176471
    $ws.Range("J13").Value = 0.3725490196078431
    $ws.Range("S13").Value = 0.0392156862745098
    $ws.Range("F15").Value = 0.04089219330855019
    $ws.Range("H15").Value = 0.1672862453531599
    $ws.Range("I15").Value = 0.02230483271375465
    $ws.Range("J15").Value = 0.3382899628252788
    $ws.Range("K15").Value = 0.05947955390334572
    $ws.Range("M15").Value = 0.02230483271375465
    $ws.Range("O15").Value = 0.09665427509293681
    $ws.Range("S15").Value = 0.2527881040892193
    $ws.Range("F16").Value = 0.01376146788990826
    $ws.Range("H16").Value = 0.2155963302752294
    $ws.Range("I16").Value = 0.04587155963302753
    $ws.Range("J16").Value = 0.3807339449541284
    $ws.Range("K16").Value = 0.1238532110091743
    $ws.Range("M16").Value = 0.02293577981651376
    $ws.Range("N16").Value = 0.004587155963302753
    $ws.Range("O16").Value = 0.04128440366972477
    $ws.Range("S16").Value = 0.1513761467889908
    $ws.Range("F17").Value = 0.01706484641638225
    $ws.Range("H17").Value = 0.2098976109215017
    $ws.Range("I17").Value = 0.08361774744027303
    $ws.Range("J17").Value = 0.4078498293515359
    $ws.Range("K17").Value = 0.1023890784982935
    $ws.Range("M17").Value = 0.01706484641638225
    $ws.Range("O17").Value = 0.06996587030716724
    $ws.Range("S17").Value = 0.09215017064846416
    $ws.Range("F18").Value = 0.01687763713080169
    $ws.Range("H18").Value = 0.2067510548523207
    $ws.Range("I18").Value = 0.03375527426160337
    $ws.Range("J18").Value = 0.3966244725738396
    $ws.Range("K18").Value = 0.109704641350211
    $ws.Range("M18").Value = 0.03375527426160337
    $ws.Range("O18").Value = 0.08016877637130802
    $ws.Range("S18").Value = 0.1223628691983122
    $ws.Range("F19").Value = 0.01337792642140468
    $ws.Range("H19").Value = 0.2615384615384616
    $ws.Range("I19").Value = 0.05752508361204013
    $ws.Range("J19").Value = 0.3551839464882943
    $ws.Range("K19").Value = 0.1311036789297659
    $ws.Range("M19").Value = 0.01806020066889632
    $ws.Range("N19").Value = 0.001337792642140468
    $ws.Range("O19").Value = 0.0588628762541806
    $ws.Range("S19").Value = 0.1030100334448161

